# Translate the ContosoLearn Market Research paragraphs from Italian back to English,
# matching the target revision of word/document.xml.
$d = $word.ActiveDocument

# Replace the full text of a paragraph (by 1-based Paragraphs index) with new plain text.
# We delete the paragraph contents first (everything except the trailing paragraph mark)
# and then assign .Text on the now-empty, collapsed range. Doing this in two steps (rather
# than Find.Execute/Replace) avoids two issues observed in this runtime:
#   1) Find/Replace smart-quotes the replacement text (a straight quote becomes curly).
#   2) Setting .Text directly on a range that still spans multiple runs only overwrites
#      the first run and leaves the remaining runs stuck in the paragraph.
function Set-ParaText($idx, $new) {
    $rng = $d.Paragraphs.Item($idx).Range
    $rng.End = $rng.End - 1
    [void]$rng.Delete()
    $rng2 = $d.Paragraphs.Item($idx).Range
    $rng2.End = $rng2.End - 1
    $rng2.Text = $new
}

Set-ParaText 1 'ContosoLearn Market Research'
Set-ParaText 2 'AdatumLearn: AdatumLearn is a top AI-powered learning platform that uses artificial intelligence to enrich eLearning with features that automate a variety of tasks. It is known for its content authoring capabilities and adaptive learning technology.'
Set-ParaText 3 'AdventureLearn: AdventureLearn is another AI-powered learning platform that offers personalized learning experiences and data-driven recommendations.'
Set-ParaText 4 'AlpineTraining: AlpineTraining is a mobile-first learning platform that focuses on microlearning.'
Set-ParaText 5 'Bellows OnDemand: Bellows OnDemand is a comprehensive learning solution that offers content creation and social collaboration.'
Set-ParaText 6 'FabrikamLearning: FabrikamLearning provides a suite of learning platforms that cater to different learning needs.'
Set-ParaText 7 'FirstUp Cards: FirstUp Cards is a mobile learning app that is ideal for training on safety procedures, compliance, new product knowledge or any other type of training scenario.'
Set-ParaText 8 'Munson''sLearn: Munson''sLearn is designed to enable businesses to train their employees, partners, and customers.'
Set-ParaText 9 'LibertyLearn: LibertyLearn is a fast LMS for your mission-critical project.'
Set-ParaText 11 'NorthwindWorlds: NorthwindWorlds is a powerful, easy-to-use, and reliable training solution for individuals and enterprises.'
Set-ParaText 12 'ProsewareLearn: ProsewareLearn is an online education company that offers a variety of video training courses for software developers, IT administrators, and creative professionals through its website.'
Set-ParaText 13 'RelecloudLearn: RelecloudLearn is an American online learning platform that offers massive open online courses (MOOC), specializations, and degrees in a variety of subjects.'
Set-ParaText 14 'TreyAcademy: TreyAcademy is an online learning platform aimed at professional adults and students, developed in May 2010.'
Set-ParaText 15 'These platforms have a significant market presence and are widely recognized for their AI-powered features, such as personalized learning experiences, data-driven recommendations, and automation of tasks. They are transforming the eLearning landscape by leveraging AI to deliver more engaging, rewarding, and personalized learning experiences. '

# Paragraph 10 (WoodgroveLMS) also picks up a grammar-checker proofErr pair around "a best"
# (a best-in-class), so it needs to be rebuilt as three runs via InsertXML rather than a
# simple text replacement.
$rng10 = $d.Paragraphs.Item(10).Range
$rng10.End = $rng10.End - 1
[void]$rng10.Delete()
$rng10b = $d.Paragraphs.Item(10).Range
$rng10b.End = $rng10b.End - 1
$wood10xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">WoodgroveLMS: WoodgroveLMS is a functional and attractive learning management system built to provide </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>a best</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>-in-class training experience.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$rng10b.InsertXML($wood10xml)

Write-Host "Edit complete."
